$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for runs(C)/balls(D)/fours(E)/sixes(F) on rows 3,5,6,7,8,9,10,12
# per the commit diff. Values are written with a leading apostrophe so Excel
# keeps them as text (matching the original numbers-stored-as-text cells)
# instead of auto-converting the numeric-looking strings to real numbers.

$ws.Range("C3").Formula = "'1"
$ws.Range("D3").Formula = "'2"
$ws.Range("E3").Formula = "'0"

$ws.Range("C5").Formula = "'15"
$ws.Range("D5").Formula = "'5"
$ws.Range("E5").Formula = "'0"

$ws.Range("C6").Formula = "'31"
$ws.Range("D6").Formula = "'21"
$ws.Range("E6").Formula = "'3"
$ws.Range("F6").Formula = "'2"

$ws.Range("C7").Formula = "'0"
$ws.Range("D7").Formula = "'3"
$ws.Range("E7").Formula = "'0"
$ws.Range("F7").Formula = "'0"

$ws.Range("C8").Formula = "'13"
$ws.Range("D8").Formula = "'14"
$ws.Range("E8").Formula = "'1"

$ws.Range("C9").Formula = "'18"
$ws.Range("D9").Formula = "'6"
$ws.Range("E9").Formula = "'1"

$ws.Range("C10").Formula = "'17"

$ws.Range("C12").Formula = "'22"
$ws.Range("D12").Formula = "'25"
$ws.Range("E12").Formula = "'1"
$ws.Range("F12").Formula = "'1"
